$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add()
$ws.Name = "Test"

$ws.Range("A2").Value = "wrap-only"
$ws.Range("A2").WrapText = $true

$ws.Range("A3").Value = "right-center"
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").HorizontalAlignment = -4152

$ws.Range("A4").Value = "center-center-wrap"
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").WrapText = $true

$ws.Range("A5").Value = "vert-center-only"
$ws.Range("A5").VerticalAlignment = -4108

Write-Output "done"
